# Generate Report for Handoff
#
# The b3b54649-5d53-4a36-bbf9-a905c1c61a4a file finished a new handoff to
# zh-cn, so its "Latest Handoff Datetime" on the zh-cn sheet is refreshed
# to reflect the newest handoff timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("E6").Value = "2016-03-24 02:44:14"
